$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 205 (new weekly price entries),
# pushing the existing rows 205-218 down to 209-222.
$ws.Rows("205:208").Insert()

# Row 205: Repollo, Copenhague, Primera — new week (2021-10-22)
$ws.Cells.Item(205, 1).Value = 11
$ws.Cells.Item(205, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(205, 3).Value = "Bíobío"
$ws.Cells.Item(205, 4).Value = 44491
$ws.Cells.Item(205, 5).Value = 8
$ws.Cells.Item(205, 6).Value = 100112006
$ws.Cells.Item(205, 7).Value = "Repollo"
$ws.Cells.Item(205, 8).Value = "Copenhague"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 1000
$ws.Cells.Item(205, 11).Value = 800
$ws.Cells.Item(205, 12).Value = 900
$ws.Cells.Item(205, 13).Value = 850
$ws.Cells.Item(205, 14).Value = "$/unidad"
$ws.Cells.Item(205, 15).Value = "Región Metropolitana"
$ws.Cells.Item(205, 16).Value = 850
$ws.Cells.Item(205, 17).Value = 1
$ws.Cells.Item(205, 18).Value = "Hortaliza"

# Row 206: Repollo, Copenhague, Segunda — new week (2021-10-22)
$ws.Cells.Item(206, 1).Value = 11
$ws.Cells.Item(206, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(206, 3).Value = "Bíobío"
$ws.Cells.Item(206, 4).Value = 44491
$ws.Cells.Item(206, 5).Value = 8
$ws.Cells.Item(206, 6).Value = 100112006
$ws.Cells.Item(206, 7).Value = "Repollo"
$ws.Cells.Item(206, 8).Value = "Copenhague"
$ws.Cells.Item(206, 9).Value = "Segunda"
$ws.Cells.Item(206, 10).Value = 500
$ws.Cells.Item(206, 11).Value = 700
$ws.Cells.Item(206, 12).Value = 700
$ws.Cells.Item(206, 13).Value = 700
$ws.Cells.Item(206, 14).Value = "$/unidad"
$ws.Cells.Item(206, 15).Value = "Región Metropolitana"
$ws.Cells.Item(206, 16).Value = 700
$ws.Cells.Item(206, 17).Value = 1
$ws.Cells.Item(206, 18).Value = "Hortaliza"

# Row 207: Repollo, Crespo record, Primera — new week (2021-10-22)
$ws.Cells.Item(207, 1).Value = 11
$ws.Cells.Item(207, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(207, 3).Value = "Bíobío"
$ws.Cells.Item(207, 4).Value = 44491
$ws.Cells.Item(207, 5).Value = 8
$ws.Cells.Item(207, 6).Value = 100112006
$ws.Cells.Item(207, 7).Value = "Repollo"
$ws.Cells.Item(207, 8).Value = "Crespo record"
$ws.Cells.Item(207, 9).Value = "Primera"
$ws.Cells.Item(207, 10).Value = 1000
$ws.Cells.Item(207, 11).Value = 800
$ws.Cells.Item(207, 12).Value = 900
$ws.Cells.Item(207, 13).Value = 850
$ws.Cells.Item(207, 14).Value = "$/unidad"
$ws.Cells.Item(207, 15).Value = "Región Metropolitana"
$ws.Cells.Item(207, 16).Value = 850
$ws.Cells.Item(207, 17).Value = 1
$ws.Cells.Item(207, 18).Value = "Hortaliza"

# Row 208: Repollo, Crespo record, Segunda — new week (2021-10-22)
$ws.Cells.Item(208, 1).Value = 11
$ws.Cells.Item(208, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(208, 3).Value = "Bíobío"
$ws.Cells.Item(208, 4).Value = 44491
$ws.Cells.Item(208, 5).Value = 8
$ws.Cells.Item(208, 6).Value = 100112006
$ws.Cells.Item(208, 7).Value = "Repollo"
$ws.Cells.Item(208, 8).Value = "Crespo record"
$ws.Cells.Item(208, 9).Value = "Segunda"
$ws.Cells.Item(208, 10).Value = 500
$ws.Cells.Item(208, 11).Value = 700
$ws.Cells.Item(208, 12).Value = 700
$ws.Cells.Item(208, 13).Value = 700
$ws.Cells.Item(208, 14).Value = "$/unidad"
$ws.Cells.Item(208, 15).Value = "Región Metropolitana"
$ws.Cells.Item(208, 16).Value = 700
$ws.Cells.Item(208, 17).Value = 1
$ws.Cells.Item(208, 18).Value = "Hortaliza"
